$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "sexo" column (F) and "edad-grupos-quinquenales" column (H) move from
# being curated as iaest-dimension to iaest-measure.
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("H2").Value = "iaest-measure:edad-grupos-quinquenales"

# Their dim/medida marker flips accordingly.
$ws.Range("F3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Their datatype switches from the dimension's skos:Concept to the
# measure's xsd:int.
$ws.Range("F4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"

# Measures have no mapping file, so the old mapping-sexo.xlsx /
# mapping-edad-grupos-quinquenales.xlsx entries disappear entirely.
$ws.Range("F5").Clear()
$ws.Range("H5").Clear()
